$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.676.68"
$ws.Range("E2").Value = "'  +0.59%  "
$ws.Range("D3").Value = "'1.927.80"
$ws.Range("E3").Value = "'  +0.19%  "
$ws.Range("E4").Value = "'  +0.58%  "
$ws.Range("D5").Value = "'326.96"
$ws.Range("E5").Value = "'  +0.31%  "
$ws.Range("D6").Value = "'1.011"
$ws.Range("E6").Value = "'  +0.56%  "
$ws.Range("D7").Value = "'0.4833"
$ws.Range("D8").Value = "'0.4061"
$ws.Range("E8").Value = "'  -0.93%  "
$ws.Range("D9").Value = "'0.08206"
$ws.Range("E9").Value = "'  +0.19%  "
$ws.Range("D10").Value = "'1.011"
$ws.Range("E10").Value = "'  -1.53%  "
$ws.Range("E11").Value = "'  -0.37%  "
$ws.Range("B12").Value = "'WrappedEther"
$ws.Range("C12").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.973.31"
$ws.Range("E12").Value = "'  +3.02%  "
$ws.Range("B13").Value = "'Polkadot"
$ws.Range("C13").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'6.079"
$ws.Range("E13").Value = "'  +0.35%  "
$ws.Range("D14").Value = "'7.297"
$ws.Range("E14").Value = "'  +0.88%  "
$ws.Range("D15").Value = "'91.42"
$ws.Range("E15").Value = "'  -0.13%  "
$ws.Range("D16").Value = "'0.06879"
$ws.Range("E16").Value = "'  +1.41%  "
$ws.Range("D17").Value = "'1.013"
$ws.Range("E17").Value = "'  +0.65%  "
$ws.Range("E18").Value = "'  -0.33%  "
$ws.Range("D19").Value = "'17.64"
$ws.Range("E19").Value = "'  -1.17%  "
$ws.Range("D20").Value = "'1.010"
$ws.Range("E20").Value = "'  +0.53%  "
$ws.Range("D21").Value = "'29.684.75"
$ws.Range("E21").Value = "'  +0.55%  "
$ws.Range("D22").Value = "'5.659"
$ws.Range("E22").Value = "'  +0.35%  "
$ws.Range("D23").Value = "'12.00"
$ws.Range("E23").Value = "'  +1.71%  "
$ws.Range("D24").Value = "'2.202"
$ws.Range("E24").Value = "'  +0.91%  "
$ws.Range("D25").Value = "'2.118.68"
$ws.Range("E25").Value = "'  -1.40%  "
$ws.Range("D26").Value = "'156.40"
$ws.Range("E26").Value = "'  -0.32%  "
$ws.Range("D27").Value = "'6.409"
$ws.Range("E27").Value = "'  -4.47%  "
$ws.Range("D28").Value = "'20.00"
$ws.Range("E28").Value = "'  -0.66%  "
$ws.Range("D29").Value = "'2.094"
$ws.Range("E29").Value = "'  -1.65%  "
$ws.Range("D30").Value = "'120.87"
$ws.Range("E30").Value = "'  +0.08%  "
$ws.Range("D31").Value = "'1.009"
$ws.Range("E31").Value = "'  -2.08%  "
$ws.Range("D32").Value = "'0.09605"
$ws.Range("E32").Value = "'  +0.20%  "
$ws.Range("D33").Value = "'5.615"
$ws.Range("E33").Value = "'  +0.97%  "
$ws.Range("B34").Value = "'ARBITRUM"
$ws.Range("C34").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.395"
$ws.Range("E34").Value = "'  -0.08%  "
$ws.Range("B35").Value = "'HuobiToken"
$ws.Range("C35").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'3.554"
$ws.Range("E35").Value = "'  -0.45%  "
$ws.Range("D36").Value = "'0.06523"
$ws.Range("E36").Value = "'  +6.02%  "
$ws.Range("D37").Value = "'0.02284"
$ws.Range("E37").Value = "'  -0.32%  "
$ws.Range("D38").Value = "'1.209"
$ws.Range("E38").Value = "'  +2.08%  "
$ws.Range("D39").Value = "'0.5937"
$ws.Range("E39").Value = "'  -0.94%  "
$ws.Range("D40").Value = "'10.76"
$ws.Range("E40").Value = "'  -0.58%  "
$ws.Range("D41").Value = "'7.876"
$ws.Range("E41").Value = "'  -2.14%  "
$ws.Range("D42").Value = "'2.556"
$ws.Range("E42").Value = "'  +3.95%  "
$ws.Range("D43").Value = "'0.1846"
$ws.Range("E43").Value = "'  -1.15%  "
$ws.Range("D44").Value = "'1.244"
$ws.Range("E44").Value = "'  -2.88%  "
$ws.Range("D45").Value = "'0.07523"
$ws.Range("E45").Value = "'  -1.41%  "
$ws.Range("D46").Value = "'12.33"
$ws.Range("D47").Value = "'0.5556"
$ws.Range("E47").Value = "'  -0.84%  "
$ws.Range("D48").Value = "'1.966"
$ws.Range("E48").Value = "'  -0.11%  "
$ws.Range("D49").Value = "'118.46"
$ws.Range("E49").Value = "'  +1.24%  "
$ws.Range("D50").Value = "'2.425"
$ws.Range("E50").Value = "'  -0.69%  "
$ws.Range("D51").Value = "'72.09"
$ws.Range("E51").Value = "'  -1.29%  "
